$wb = $excel.ActiveWorkbook

# --- Config sheet ("sheet1" in the package): add a "Cycle" column ---
$wsConfig = $wb.Worksheets.Item("Config")
$wsConfig.Range("E1").Value = "Cycle"
$wsConfig.Range("E2").Value = 50

# --- UUT_WES7P-64 sheet ("sheet3" in the package): move the selection ---
# Select on this sheet first (without leaving it as the active/tab-selected
# sheet) so only its <selection> changes while "Config" stays the tab shown.
$wsPlan = $wb.Worksheets.Item("UUT_WES7P-64")
$wsPlan.Range("F26").Select()

# --- Re-activate Config and move its selection last so it remains the ---
# --- active (tabSelected) sheet with its own updated selection.        ---
$wsConfig.Activate()
$wsConfig.Range("B21").Select()
